$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (this also updates the defined name's sheet reference
# and the workbook title automatically).
$ws.Name = "20160930 OIH"

# --- Update existing FIELD/VALUE rows for the new OIH trade ---
$ws.Range("B2").Value = "20160930 +OIH-161021C30.00"
$ws.Range("B3").Value = "E:\\Datos\\bolsa\\cuenta personal\\analisis de valores\\Trades activos\\Scanning\\20160930"
$ws.Range("B4").Value = "call"
$ws.Range("B5").Value = 30
$ws.Range("B9").Value = 0.42
$ws.Range("B10").Value = 29.05
$ws.Range("B16").Value = 36
$ws.Range("B17").Value = "OIH"
$ws.Range("B18").Value = 0.3767

# --- Fill in the previously-empty bid/ask (write leg) rows 19-30 ---
$ws.Range("B19").Value = "20160930 +OIH-161021C29.00"
$ws.Range("B20").Value = "E:\\Datos\\bolsa\\cuenta personal\\analisis de valores\\Trades activos\\Scanning\\20160930"
$ws.Range("B21").Value = 29
$ws.Range("B22").Value = -0.86
$ws.Range("B23").Value = 29.035
$ws.Range("B24").Value = 2016
$ws.Range("B25").Value = 9
$ws.Range("B26").Value = 30
$ws.Range("B27").Value = 10
$ws.Range("B28").Value = 24
$ws.Range("B29").Value = 0
$ws.Range("B30").Value = 0.3767

# Move the active selection to reflect where the edits finished.
[void]$ws.Range("B23").Select()
